$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column before DY (24-nov) ---
$ws = $wb.Worksheets.Item("Prix Spot")
$ws.Range("DY1").EntireColumn.Insert()
$ws.Range("DY1").Value = "24-nov"
$ws.Range("DY2:DY25").Value = "-"

# --- Sheet "Gaz": append two new daily rows ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A158").NumberFormat = "@"
$wsGaz.Range("A158").Value = "2025-11-22"
$wsGaz.Range("A158").ClearFormats()
$wsGaz.Range("B158").Value = 29.2

$wsGaz.Range("A159").NumberFormat = "@"
$wsGaz.Range("A159").Value = "2025-11-23"
$wsGaz.Range("A159").ClearFormats()
$wsGaz.Range("B159").Value = 29.2

# --- Sheet "CO2": append two new daily rows ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A158").NumberFormat = "@"
$wsCo2.Range("A158").Value = "2025-11-22"
$wsCo2.Range("A158").ClearFormats()
$wsCo2.Range("B158").Value = 80.28

$wsCo2.Range("A159").NumberFormat = "@"
$wsCo2.Range("A159").Value = "2025-11-23"
$wsCo2.Range("A159").ClearFormats()
$wsCo2.Range("B159").Value = 80.28
